# Insert two new weekly price rows for "Hortaliza, Agrícola del Norte S.A.
# de Arica - Locoto" immediately after the current row 110 (pushing the
# existing rows 111-139 down to 113-141), then populate the two new rows
# with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 111 (row 111 becomes blank, old 111
# shifts to 112, then inserting again at 111 pushes everything down by 2
# total, so old row 111 ends up at 113).
$ws.Rows.Item(111).Insert()
$ws.Rows.Item(111).Insert()

# New row 111: Locoto, Primera
$ws.Range("A111").Value = 1
$ws.Range("B111").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C111").Value = "Arica y Parinacota"
$ws.Range("D111").Value = "2022-10-07"
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 100112042
$ws.Range("G111").Value = "Locoto"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 130
$ws.Range("K111").Value = 20000
$ws.Range("L111").Value = 21000
$ws.Range("M111").Value = 20500
$ws.Range("N111").Value = "$/caja 20 kilos"
$ws.Range("O111").Value = "Región de Arica y Parinacota"
$ws.Range("P111").Value = 1025
$ws.Range("Q111").Value = 20
$ws.Range("R111").Value = "Hortaliza"

# New row 112: Locoto, Segunda
$ws.Range("A112").Value = 1
$ws.Range("B112").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C112").Value = "Arica y Parinacota"
$ws.Range("D112").Value = "2022-10-07"
$ws.Range("E112").Value = 15
$ws.Range("F112").Value = 100112042
$ws.Range("G112").Value = "Locoto"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Segunda"
$ws.Range("J112").Value = 150
$ws.Range("K112").Value = 18000
$ws.Range("L112").Value = 19000
$ws.Range("M112").Value = 18500
$ws.Range("N112").Value = "$/caja 20 kilos"
$ws.Range("O112").Value = "Región de Arica y Parinacota"
$ws.Range("P112").Value = 925
$ws.Range("Q112").Value = 20
$ws.Range("R112").Value = "Hortaliza"
